$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.687.59"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.883.28"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'331.20"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4733"
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("D8").Value = "'0.3978"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "'47.55"
$ws.Range("E9").Value = "  -4.97%  "
$ws.Range("D10").Value = "'0.08037"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("D11").Value = "'1.025"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").Value = "'21.83"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.879.25"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "'5.972"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "'7.181"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'87.00"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "'0.00001041"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "'0.06612"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'17.22"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "27.705.34"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "'5.511"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").Value = "2.088.95"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").Value = "'155.63"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").Value = "'20.28"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").Value = "'2.095"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").Value = "'5.565"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").Value = "'122.37"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'0.9675"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "'0.09554"
$ws.Range("D34").Value = "'1.471"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").Value = "'3.629"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "'5.308"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("D37").Value = "'0.06128"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'0.02255"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "'1.226"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").Value = "'8.142"
$ws.Range("E40").Value = "  -5.43%  "
$ws.Range("D41").Value = "'0.6007"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'0.1900"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'10.31"
$ws.Range("E44").Value = "  -4.18%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5699"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.245"
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("D47").Value = "'12.28"
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("D48").Value = "'3.405"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "'1.934"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").Value = "'0.06824"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "'110.91"
$ws.Range("E51").Value = "  +0.41%  "
